$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Existing rows 8 & 9 (OwnProfileCommentsLikeTest / OthersProfileCommentsLikeTest)
# move from Runmode "Y" / Result "PASS" to Runmode "N" / Result "SKIP"
$ws.Range("C8").Value = "N"
$ws.Range("D8").Value = "SKIP"
$ws.Range("C9").Value = "N"
$ws.Range("D9").Value = "SKIP"

# New test case row: Profile Interest and Skills Update
$ws.Range("A10").Value = "ProfileInterestSkillsUpdateTest"
$ws.Range("B10").Value = "To verify User can able to update his Own Profile Interests and Skills"
$ws.Range("C10").Value = "Y"
$ws.Range("D10").Value = "PASS"

# Match the formatting used by the row above it
$ws.Range("A9:D9").Copy()
$ws.Range("A10:D10").PasteSpecial(-4122)

$ws.Range("B12").Select()
